$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Add Devices Loop A")
$ws2 = $wb.Worksheets.Item("Add_Devices_LoopB_PFI")
$ws3 = $wb.Worksheets.Item("Add_Devices_LoopB_FIM")

$sheets = @($ws1, $ws2, $ws3)

foreach ($ws in $sheets) {
    # Copy formatting from existing header/data cells in column I
    $ws.Range("I7").Copy() | Out-Null
    $ws.Range("I1").PasteSpecial(-4122) | Out-Null

    $ws.Range("I8").Copy() | Out-Null
    $ws.Range("I2:I3").PasteSpecial(-4122) | Out-Null

    $ws.Range("I1").Value = "DC Unit Loading Details Name"
    $ws.Range("I2").Value = "Current (DC Units)"
    $ws.Range("I3").Value = "Current (worst case)"
}

$excel.CutCopyMode = 0

# Update selections / scroll position for each sheet.
# Visit sheet3, then sheet2, and finally sheet1 so that "Add Devices Loop A"
# ends up as the active / tab-selected sheet, matching the final state in Excel.
$ws3.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 4
$ws3.Range("I1:I3").Select() | Out-Null

$ws2.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 4
$ws2.Range("I1:I3").Select() | Out-Null

# "Add Devices Loop A" is the final active sheet
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 4
$ws1.Range("I1:I3").Select() | Out-Null
